$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 451, shifting existing rows 451:507 down to 452:508
$ws.Rows.Item(451).Insert()

# Populate the newly inserted row 451 with the new data record
$ws.Range("A451").Value2 = 4
$ws.Range("B451").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C451").Value2 = "Los Lagos"
$ws.Range("D451").Value2 = 45142
$ws.Range("E451").Value2 = 10
$ws.Range("F451").Value2 = 100112040
$ws.Range("G451").Value2 = "Cilantro"
$ws.Range("H451").Value2 = "Sin especificar"
$ws.Range("I451").Value2 = "Primera"
$ws.Range("J451").Value2 = 180
$ws.Range("K451").Value2 = 11000
$ws.Range("L451").Value2 = 12000
$ws.Range("M451").Value2 = 11500
$ws.Range("N451").Value2 = "$/caja 36 atados"
$ws.Range("O451").Value2 = "Región Metropolitana"
$ws.Range("P451").Value2 = 319
$ws.Range("Q451").Value2 = 36
$ws.Range("R451").Value2 = "Hortaliza"
